$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 316.5
$ws.Range("J4").Value = 999
$ws.Range("L4").Value = 999
$ws.Range("N4").Value = -1227
$ws.Range("H11").Value = 300.41177
$ws.Range("I11").Value = 300.41177
$ws.Range("K11").Value = 300.41177
$ws.Range("M11").Value = -160.41177
$ws.Range("H19").Value = 2239.3333
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 2239.3333
$ws.Range("K19").Value = 0
$ws.Range("L19").ClearContents()
$ws.Range("M19").Value = 2239.3333
$ws.Range("N19").Value = -2589.3333
$ws.Range("H32").Value = 6944.2
$ws.Range("J32").Value = 6944.2
$ws.Range("L32").Value = 6944.2
$ws.Range("N32").Value = -7596.2
$ws.Range("H33").Value = 25000612
$ws.Range("I33").Value = 25000612
$ws.Range("K33").Value = 25000612
$ws.Range("M33").Value = -25000383
$ws.Range("H40").Value = 2919.2
$ws.Range("J40").Value = 2924.25
$ws.Range("L40").Value = 2924.25
$ws.Range("N40").Value = -3274.25
$ws.Range("H109").Value = 44921
$ws.Range("J109").Value = 44921
$ws.Range("L109").Value = 44921
$ws.Range("N109").Value = -47695
$ws.Range("H138").Value = 8067870
$ws.Range("J138").Value = 12825150
$ws.Range("L138").Value = 38475450
$ws.Range("N138").Value = -38485730

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 66672840
$ws.Range("I2").Value = 76928664
$ws.Range("J2").Value = 9995
$ws.Range("K2").Value = 76928664
$ws.Range("L2").Value = 9995
$ws.Range("M2").Value = -76928551
$ws.Range("N2").Value = -10221
$ws.Range("H61").Value = 3899.6
$ws.Range("J61").Value = 4531.6924
$ws.Range("L61").Value = 4531.6924
$ws.Range("N61").Value = -4955.6924
$ws.Range("H74").Value = 29672.691
$ws.Range("I74").Value = 29393.305
$ws.Range("J74").Value = 31349
$ws.Range("K74").Value = 29393.305
$ws.Range("L74").Value = 31349
$ws.Range("M74").Value = -28519.305
$ws.Range("N74").Value = -33097
$ws.Range("H77").Value = 29672.691
$ws.Range("I77").Value = 29393.305
$ws.Range("J77").Value = 31349
$ws.Range("K77").Value = 146966.525
$ws.Range("L77").Value = 156745
$ws.Range("M77").Value = -142598.525
$ws.Range("N77").Value = -165481
$ws.Range("H116").Value = 66672840
$ws.Range("I116").Value = 76928664
$ws.Range("J116").Value = 9995
$ws.Range("K116").Value = 76928664
$ws.Range("L116").Value = 9995
$ws.Range("M116").Value = -76926370
$ws.Range("N116").Value = -14583
$ws.Range("H132").Value = 2720.9023
$ws.Range("I132").Value = 2458.6572
$ws.Range("J132").Value = 4250.6665
$ws.Range("K132").Value = 7375.971600000001
$ws.Range("L132").Value = 12751.9995
$ws.Range("M132").Value = -4845.971600000001
$ws.Range("N132").Value = -17811.9995
$ws.Range("H136").Value = 3899.6
$ws.Range("J136").Value = 4531.6924
$ws.Range("L136").Value = 13595.0772
$ws.Range("N136").Value = -18695.0772

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 66672840
$ws.Range("I3").Value = 76928664
$ws.Range("J3").Value = 9995
$ws.Range("K3").Value = 76928664
$ws.Range("L3").Value = 9995
$ws.Range("M3").Value = -76928550
$ws.Range("N3").Value = -10223
$ws.Range("H118").Value = 99000
$ws.Range("J118").Value = 99000
$ws.Range("L118").Value = 99000
$ws.Range("N118").Value = -102314
$ws.Range("H134").Value = 2145.1667
$ws.Range("I134").Value = 2145.1667
$ws.Range("K134").Value = 6435.500100000001
$ws.Range("M134").Value = -3900.500100000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20300.35
$ws.Range("I31").Value = 30710.334
$ws.Range("J31").Value = 4685.375
$ws.Range("K31").Value = 30710.334
$ws.Range("L31").Value = 4685.375
$ws.Range("M31").Value = -30415.334
$ws.Range("N31").Value = -5275.375
$ws.Range("H34").Value = 20300.35
$ws.Range("I34").Value = 30710.334
$ws.Range("J34").Value = 4685.375
$ws.Range("K34").Value = 30710.334
$ws.Range("L34").Value = 4685.375
$ws.Range("M34").Value = -30508.334
$ws.Range("N34").Value = -5089.375
$ws.Range("H132").Value = 2652.3809
$ws.Range("I132").Value = 2418
$ws.Range("J132").Value = 3402.4
$ws.Range("K132").Value = 7254
$ws.Range("L132").Value = 10207.2
$ws.Range("M132").Value = -4724
$ws.Range("N132").Value = -15267.2

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 6127
$ws.Range("I64").Value = 4337.3335
$ws.Range("J64").Value = 7916.6665
$ws.Range("K64").Value = 13012.0005
$ws.Range("L64").Value = 23749.9995
$ws.Range("M64").Value = -12742.0005
$ws.Range("N64").Value = -24289.9995
$ws.Range("H67").Value = 6127
$ws.Range("I67").Value = 4337.3335
$ws.Range("J67").Value = 7916.6665
$ws.Range("K67").Value = 13012.0005
$ws.Range("L67").Value = 23749.9995
$ws.Range("M67").Value = -12076.0005
$ws.Range("N67").Value = -25621.9995
$ws.Range("H107").Value = 914
$ws.Range("I107").Value = 696.8
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 2090.4
$ws.Range("L107").Value = 6000
$ws.Range("M107").Value = -170.3999999999996
$ws.Range("N107").Value = -9840

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 76925230
$ws.Range("I102").Value = 2334.25
$ws.Range("J102").Value = 1000000000
$ws.Range("K102").Value = 2334.25
$ws.Range("L102").Value = 1000000000
$ws.Range("M102").Value = -712.25
$ws.Range("N102").Value = -1000003244
$ws.Range("H122").Value = 3990
$ws.Range("I122").Value = 3990
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 11970
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -9520
$ws.Range("H126").Value = 3335.8667
$ws.Range("I126").Value = 3003.3333
$ws.Range("J126").Value = 4666
$ws.Range("K126").Value = 9009.999899999999
$ws.Range("L126").Value = 13998
$ws.Range("M126").Value = -6539.999899999999
$ws.Range("N126").Value = -18938
$ws.Range("H128").Value = 97500
$ws.Range("J128").Value = 97500
$ws.Range("L128").Value = 97500
$ws.Range("N128").Value = -107460
$ws.Range("H130").Value = 61999
$ws.Range("J130").Value = 61999
$ws.Range("L130").Value = 61999
$ws.Range("N130").Value = -72039
$ws.Range("H132").Value = 4902.5
$ws.Range("I132").Value = 4985.0586
$ws.Range("J132").Value = 3499
$ws.Range("K132").Value = 14955.1758
$ws.Range("L132").Value = 10497
$ws.Range("M132").Value = -12425.1758
$ws.Range("N132").Value = -15557

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6920.8335
$ws.Range("I7").Value = 6472
$ws.Range("K7").Value = 6472
$ws.Range("M7").Value = -6360
$ws.Range("H100").Value = 10500.333
$ws.Range("J100").Value = 14124.75
$ws.Range("L100").Value = 14124.75
$ws.Range("N100").Value = -15206.75
$ws.Range("H126").Value = 6920.8335
$ws.Range("I126").Value = 6472
$ws.Range("K126").Value = 19416
$ws.Range("M126").Value = -16946

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 74142
$ws.Range("J41").Value = 74771.164
$ws.Range("L41").Value = 74771.164
$ws.Range("N41").Value = -75551.164
$ws.Range("H96").Value = 1424.5
$ws.Range("I96").Value = 1203.9166
$ws.Range("K96").Value = 1203.9166
$ws.Range("M96").Value = 169.0834
$ws.Range("H122").Value = 1532.0968
$ws.Range("I122").Value = 1303.3928
$ws.Range("K122").Value = 3910.1784
$ws.Range("M122").Value = -1460.1784
$ws.Range("H132").Value = 1423.3036
$ws.Range("I132").Value = 1421.1459
$ws.Range("K132").Value = 4263.4377
$ws.Range("M132").Value = -1733.4377
$ws.Range("H136").Value = 1518.0143
$ws.Range("I136").Value = 1313.638
$ws.Range("K136").Value = 3940.914
$ws.Range("M136").Value = -1390.914
